$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "28.194.37"
$ws.Cells.Item(2, 5).Value = "  -3.06%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.915.28"
$ws.Cells.Item(3, 5).Value = "  -3.69%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "327.97"
$ws.Cells.Item(5, 5).Value = "  -0.68%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(6, 5).Value = "  -1.03%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4679"
$ws.Cells.Item(7, 5).Value = "  -6.00%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.4015"
$ws.Cells.Item(8, 5).Value = "  -4.37%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "53.21"
$ws.Cells.Item(9, 5).Value = "  -1.97%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.08401"
$ws.Cells.Item(10, 5).Value = "  -10.86%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.044"
$ws.Cells.Item(11, 5).Value = "  -4.74%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -4.19%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.893.61"
$ws.Cells.Item(13, 5).Value = "  -5.86%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.412"
$ws.Cells.Item(14, 5).Value = "  -7.11%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "6.068"
$ws.Cells.Item(15, 5).Value = "  -5.91%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "1.001"
$ws.Cells.Item(16, 5).Value = "  -1.19%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "89.58"
$ws.Cells.Item(17, 5).Value = "  -3.24%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.00001062"
$ws.Cells.Item(18, 5).Value = "  -4.72%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06590"
$ws.Cells.Item(19, 5).Value = "  -2.19%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "18.00"
$ws.Cells.Item(20, 5).Value = "  -7.74%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.9998"
$ws.Cells.Item(21, 5).Value = "  -1.05%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.713"
$ws.Cells.Item(22, 5).Value = "  -4.44%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "28.185.95"
$ws.Cells.Item(23, 5).Value = "  -3.15%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.25"
$ws.Cells.Item(24, 5).Value = "  -6.35%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.280"
$ws.Cells.Item(25, 5).Value = "  -0.43%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "2.118.41"
$ws.Cells.Item(26, 5).Value = "  -5.25%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "153.26"
$ws.Cells.Item(27, 5).Value = "  -2.18%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -4.01%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.726"
$ws.Cells.Item(29, 5).Value = "  -8.74%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "LidoDAOToken"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.119"
$ws.Cells.Item(30, 5).Value = "  -6.48%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "123.40"
$ws.Cells.Item(31, 5).Value = "  -3.21%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.9733"
$ws.Cells.Item(32, 5).Value = "  -7.26%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.09608"
$ws.Cells.Item(33, 5).Value = "  -2.43%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.449"
$ws.Cells.Item(34, 5).Value = "  -5.66%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.631"
$ws.Cells.Item(35, 5).Value = "  -3.14%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "Filecoin"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.542"
$ws.Cells.Item(36, 5).Value = "  -4.77%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.02303"
$ws.Cells.Item(37, 5).Value = "  -5.20%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "FraxShare"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "8.762"
$ws.Cells.Item(38, 5).Value = "  -3.40%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.262"
$ws.Cells.Item(39, 5).Value = "  -4.41%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.06170"
$ws.Cells.Item(40, 5).Value = "  -3.66%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.6130"
$ws.Cells.Item(41, 5).Value = "  -5.52%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "11.03"
$ws.Cells.Item(42, 5).Value = "  -4.07%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.9997"
$ws.Cells.Item(43, 5).Value = "  -1.02%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.1902"
$ws.Cells.Item(44, 5).Value = "  -4.62%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.305"
$ws.Cells.Item(45, 5).Value = "  -3.67%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "12.83"
$ws.Cells.Item(46, 5).Value = "  -4.44%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5841"
$ws.Cells.Item(47, 5).Value = "  -6.20%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.021"
$ws.Cells.Item(48, 5).Value = "  -7.56%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.434"
$ws.Cells.Item(49, 5).Value = "  -1.77%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.06881"
$ws.Cells.Item(50, 5).Value = "  -1.63%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.00000000309"
$ws.Cells.Item(51, 5).Value = "  -9.23%  "
